# Update odds values on Sheet1 for rows 2 and 3 per FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("G2").Value = 1.6
$ws.Range("H2").Value = 3.6
$ws.Range("I2").Value = 6.25
$ws.Range("J2").Value = 2.25
$ws.Range("L2").Value = 6.5
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.67
$ws.Range("Z2").Value = 11
$ws.Range("AD2").Value = 7
$ws.Range("AF2").Value = 67
$ws.Range("AG2").Value = 13
$ws.Range("AH2").Value = 29
$ws.Range("AJ2").Value = 67
$ws.Range("AM2").Value = 3.4
$ws.Range("AN2").Value = 8.5
$ws.Range("AP2").Value = 29
$ws.Range("AT2").Value = 9.5
$ws.Range("AU2").Value = 67
$ws.Range("AV2").Value = 7
$ws.Range("AW2").Value = 34
$ws.Range("AY2").Value = 126
$ws.Range("AZ2").Value = 151

# Row 3 changes
$ws.Range("G3").Value = 1.55
$ws.Range("H3").Value = 3.7
$ws.Range("I3").Value = 7
$ws.Range("Y3").Value = 8.5
$ws.Range("AE3").Value = 21
$ws.Range("AG3").Value = 15
$ws.Range("AH3").Value = 34
$ws.Range("AY3").Value = 151
